$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay TEXT (as in the source workbook).
# A leading apostrophe forces Excel to store the literal text (quote-prefixed); resetting
# the style back to "Normal" afterwards drops the quote-prefix style flag again so the
# cell ends up with plain text and the original (default) cell style, matching the source.

$ws.Range("D2").Value = "'63.898.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "'2.745.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.40%  "

$ws.Range("D5").Value = "'576.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("D6").Value = "'157.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.89%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("D11").Value = "'0.383"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("E12").Value = "  -16.36%  "

$ws.Range("D13").Value = "'3.232.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.45%  "

$ws.Range("D14").Value = "'26.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").Value = "'63.713.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").Value = "'0.0000149"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").Value = "'2.753.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'12.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").Value = "'356.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").Value = "'6.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "

$ws.Range("D22").Value = "'0.540"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'65.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").Value = "'0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("D27").Value = "'8.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "

$ws.Range("D28").Value = "'0.0₃0893"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.32%  "

$ws.Range("D29").Value = "'1.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.45%  "

$ws.Range("D30").Value = "'6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.12%  "

$ws.Range("D31").Value = "'169.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("D32").Value = "'1.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("D33").Value = "'20.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "

$ws.Range("D34").Value = "'4.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "'1.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.41%  "

$ws.Range("D37").Value = "'1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").Value = "'0.978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.62%  "

$ws.Range("D39").Value = "'6.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.48%  "

$ws.Range("D40").Value = "'4.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("D41").Value = "'324.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.15%  "

$ws.Range("D42").Value = "'39.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.01%  "

$ws.Range("D43").Value = "'21.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.05%  "

$ws.Range("D44").Value = "'0.0588"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.55%  "

$ws.Range("D45").Value = "'21.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").Value = "'0.0253"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("D47").Value = "'135.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "

$ws.Range("D48").Value = "'0.626"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "

$ws.Range("E49").Value = "  -0.26%  "

$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").Value = "'11.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.58%  "
